# Update cryptocurrency price (D) and volume-change (E) columns
# with freshly scraped values, preserving the original text
# formatting (values are stored as text, not numbers).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.591.75"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.49%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.861.65"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.87%  "

$ws.Range("E4").Value = "  +0.73%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "334.08"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.24%  "

$ws.Range("E6").Value = "  +0.73%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4682"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.55%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3904"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.30%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "45.76"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -4.36%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07978"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.08%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.001"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.62%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "21.67"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.52%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.877.05"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.26%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.989"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.25%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.237"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.58%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.013"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.95%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "88.14"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.09%  "

$ws.Range("E18").Value = "  +0.73%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.00001043"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.17%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "16.99"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.36%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.011"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.73%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "27.587.66"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.51%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.448"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.52%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.85"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.53%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.307"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.12%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.093.45"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.03%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "159.14"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.04%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.76"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.23%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.125"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.01%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.399"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.22%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "121.49"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.39%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9741"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.97%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09459"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.52%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.614"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.52%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.281"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.46%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.330"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -7.94%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06042"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.42%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02223"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.84%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.193"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.77%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.244"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.50%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.010"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.72%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5924"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.56%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1879"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.16%  "

$ws.Range("E44").Value = "  -1.02%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.254"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.34%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5627"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.29%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "12.18"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.28%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.917"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.83%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.247"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.32%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06768"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.14%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "112.52"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.94%  "
